$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 2.2.0-ballot -> 2.1.0
$meta.Range("B3").Value = "2.1.0"

# Date: 2025-12-19T08:32:44+00:00 -> 2025-12-19T08:44:55+00:00
$meta.Range("B8").Value = "2025-12-19T08:44:55+00:00"

# Base Definition: drop the "|4.0.1" version suffix
$meta.Range("B18").Value = "http://hl7.org/fhir/StructureDefinition/Extension"

# --- Elements sheet updates ---
$elem = $wb.Worksheets.Item("Elements")

# Extension.value[x] Type(s) cell: drop the "|2.2.0-ballot" version suffix
$elem.Range("K6").Value = "Reference(https://interop.esante.gouv.fr/ig/fhir/tddui/StructureDefinition/tddui-careplan-projet-personnalise)
"

# Column K width: 95.74609375 -> 86.23828125
# (the COM width setter snaps to Excel's internal pixel grid, same as
# interactively typing a width in the UI would; 85.33 is the input that
# lands on the closest representable grid value to the target)
$elem.Columns.Item(11).ColumnWidth = 85.33
